# Rename "Second Sheet" to "Copy Sheet" and reset its A1 cell to match
# "First Sheet"'s A1 content ("Hello World!") instead of the old "New Sheet" text.

$wb = $excel.ActiveWorkbook

$wsFirst = $wb.Worksheets.Item("First Sheet")
$wsSecond = $wb.Worksheets.Item("Second Sheet")

# Update the cell content first (while we can still reference the sheet by its old name).
$wsSecond.Range("A1").Value = $wsFirst.Range("A1").Value2

# Rename the sheet.
$wsSecond.Name = "Copy Sheet"
